$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set a cell value as literal text, preserving default (unstyled) formatting
# by temporarily forcing a text number format while assigning the value, then
# reverting the style back to Normal so no explicit style index is left behind.
function Set-TextValue($cell, $value) {
    $r = $ws.Range($cell)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

Set-TextValue "D2" "64.511.83"
Set-TextValue "E2" "  +1.98%  "
Set-TextValue "D3" "2.664.08"
Set-TextValue "E3" "  +0.59%  "
Set-TextValue "E4" "  +0.07%  "
Set-TextValue "D5" "608.24"
Set-TextValue "E5" "  +0.36%  "
Set-TextValue "D6" "152.59"
Set-TextValue "E6" "  +5.67%  "
Set-TextValue "E7" "  +0.05%  "
Set-TextValue "D8" "0.592"
Set-TextValue "E8" "  +1.05%  "
Set-TextValue "E9" "  +1.98%  "
Set-TextValue "E10" "  +6.79%  "
Set-TextValue "D11" "5.61"
Set-TextValue "E11" "  -0.33%  "
Set-TextValue "E12" "  -0.71%  "
Set-TextValue "D13" "28.21"
Set-TextValue "E13" "  +3.29%  "
Set-TextValue "D14" "3.142.39"
Set-TextValue "E14" "  +0.71%  "
Set-TextValue "D15" "64.361.19"
Set-TextValue "E15" "  +2.01%  "
Set-TextValue "E16" "  +2.74%  "
Set-TextValue "D17" "2.666.65"
Set-TextValue "E17" "  +0.62%  "
Set-TextValue "D18" "12.20"
Set-TextValue "E18" "  +7.35%  "
Set-TextValue "D19" "4.66"
Set-TextValue "E19" "  +5.04%  "
Set-TextValue "D20" "348.12"
Set-TextValue "E20" "  +1.58%  "
Set-TextValue "D21" "6.96"
Set-TextValue "E21" "  +1.76%  "
Set-TextValue "D22" "0.999"
Set-TextValue "E22" "  -0.04%  "
Set-TextValue "E23" "  +0.17%  "
Set-TextValue "D24" "66.80"
Set-TextValue "E24" "  -0.46%  "
Set-TextValue "D25" "1.74"
Set-TextValue "E25" "  +13.00%  "
$ws.Range("B26").Value = "InternetComputer(DFINITY)"
$ws.Range("C26").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D26" "9.40"
Set-TextValue "E26" "  +8.72%  "
$ws.Range("B27").Value = "Fetch.AI"
$ws.Range("C27").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue "D27" "1.71"
Set-TextValue "E27" "  +4.21%  "
Set-TextValue "D28" "8.24"
Set-TextValue "E28" "  +4.22%  "
Set-TextValue "D29" "558.96"
Set-TextValue "E29" "  +2.00%  "
Set-TextValue "D30" "0.165"
Set-TextValue "E30" "  +0.25%  "
Set-TextValue "D31" "0.999"
Set-TextValue "E31" "  -0.17%  "
Set-TextValue "D32" "2.06"
Set-TextValue "E32" "  +1.39%  "
Set-TextValue "D33" "0.0₃0868"
Set-TextValue "E33" "  +7.20%  "
Set-TextValue "E34" "  -0.87%  "
Set-TextValue "D35" "5.38"
Set-TextValue "E35" "  +5.71%  "
Set-TextValue "D36" "168.62"
Set-TextValue "E36" "  -2.27%  "
Set-TextValue "E37" "  +0.86%  "
Set-TextValue "D38" "1.98"
Set-TextValue "E38" "  +7.91%  "
Set-TextValue "D39" "1.00"
Set-TextValue "E39" "  +0.04%  "
Set-TextValue "D40" "19.44"
Set-TextValue "E40" "  +1.75%  "
Set-TextValue "E41" "  +0.04%  "
Set-TextValue "D42" "167.52"
Set-TextValue "E42" "  -2.50%  "
Set-TextValue "D43" "40.36"
Set-TextValue "E43" "  +0.55%  "
Set-TextValue "D44" "3.87"
Set-TextValue "E44" "  +3.19%  "
Set-TextValue "D45" "0.0579"
Set-TextValue "E45" "  +1.28%  "
Set-TextValue "D46" "22.14"
Set-TextValue "E46" "  -0.79%  "
Set-TextValue "D47" "0.632"
Set-TextValue "E47" "  +0.02%  "
Set-TextValue "D48" "2.01"
Set-TextValue "E48" "  +15.18%  "
Set-TextValue "E49" "  +3.15%  "
Set-TextValue "D50" "0.0967"
Set-TextValue "E50" "  +0.78%  "
Set-TextValue "D51" "19.14"
Set-TextValue "E51" "  +1.83%  "
